# Update cached regression-result figures in the fc_robustness sheet.
#
# These cells hold formulas that pull their cached display text from an
# external workbook link ([1]fc_robustness!...). The COM object model has
# no supported way to poke a formula cell's cached value without a live,
# resolvable external source, so we write the refreshed figures directly
# as the new cell contents (mirroring what Excel itself bakes into the
# cell after "Edit Links > Update Values" is saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fc_robustness")

$updates = @{
    "B4"  = "-202.3***"
    "C4"  = "-296.7***"
    "D4"  = "-205.9***"
    "E4"  = "-95.6***"
    "F4"  = "-143.1*"

    "B5"  = "(48.2)"
    "C5"  = "(83.4)"
    "D5"  = "(49.1)"

    "B6"  = "-40.1"
    "C6"  = "-57.9"
    "D6"  = "-33.8"
    "E6"  = "-33.2"
    "F6"  = "-26.8"

    "E7"  = "(39.0)"
    "F7"  = "(74.2)"

    "B11" = "941.1"
    "C11" = "1387.7"
    "D11" = "1024.8"
    "E11" = "478.6"
    "F11" = "925.4"

    "E15" = "-0.058***"
    "F15" = "-0.093**"

    "B17" = "-0.0077"
    "C17" = "-0.050"
    "D17" = "-0.0021"
    "E17" = "-0.030*"
    "F17" = "-0.041"

    "C18" = "(0.044)"

    "D21" = "0.026"
    "D22" = "0.71"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
